$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: the existing "[4]张煜..." reference paragraph (last paragraph
# in the document) currently has its text split across three runs
# ("[", "4", "]张煜...."). The edit merges them into a single run
# while leaving the paragraph mark (and its paraId/rsid attributes)
# untouched.
# ---------------------------------------------------------------------
$n = $d.Paragraphs.Count
$p4 = $d.Paragraphs.Item($n)
$r4 = $p4.Range
$r4.MoveEnd(1, -1)
$origText = $r4.Text

# Force a real text-layout change so the run list is rebuilt as a
# single run (re-assigning the exact same text is treated as a no-op
# and leaves the original run split untouched).
$r4.Text = $origText + "#"
$r4b = $p4.Range
$r4b.MoveEnd(1, -1)
$r4b.Text = $origText

# ---------------------------------------------------------------------
# Step 2: append a brand-new paragraph after it (before the sectPr)
# containing the "[5] 陆雄文.管理学大辞典..." reference. This mirrors
# what Word produces when a user types the text in right after the
# previous paragraph: several runs, some carrying an explicit
# <w:rFonts w:hint="eastAsia"/> because of IME input.
# ---------------------------------------------------------------------
$tail = $d.Content
$tail.Collapse(0)

$newParaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/newpara.xml" pkg:contentType="text/xml">
    <pkg:xmlData>
      <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:pPr>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>[</w:t>
        </w:r>
        <w:r>
          <w:t>5]</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t xml:space="preserve"> </w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:rFonts w:hint="eastAsia"/>
          </w:rPr>
          <w:t>陆雄文</w:t>
        </w:r>
        <w:r>
          <w:t>.管理学大辞典:[M].上海:上海世纪出版股份有限公司上海辞书出版社,</w:t>
        </w:r>
      </w:p>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$tail.InsertXML($newParaXml)
